# Apply updated crypto price/volume figures per diff (columns D and E, rows 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.754.37"
$ws.Range("E2").Value = "  +0.28%  "

$ws.Range("D3").Value = "1.906.61"
$ws.Range("E3").Value = "  +0.64%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9979"
$ws.Range("E4").Value = "  -0.35%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.58"

$ws.Range("E6").Value = "  -0.29%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5228"
$ws.Range("E7").Value = "  +7.53%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3785"
$ws.Range("E8").Value = "  -0.28%  "

$ws.Range("E9").Value = "  -1.31%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.36"
$ws.Range("E10").Value = "  +3.93%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9024"
$ws.Range("E11").Value = "  -1.56%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07634"
$ws.Range("E12").Value = "  -0.84%  "

$ws.Range("D13").Value = "1.861.69"
$ws.Range("E13").Value = "  -1.52%  "

$ws.Range("E14").Value = "  -0.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.18"
$ws.Range("E15").Value = "  +1.27%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9980"
$ws.Range("E16").Value = "  -0.43%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008703"
$ws.Range("E17").Value = "  -1.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9986"
$ws.Range("E18").Value = "  -0.25%  "

$ws.Range("D19").Value = "27.777.12"
$ws.Range("E19").Value = "  +0.22%  "

$ws.Range("E20").Value = "  +0.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.147"
$ws.Range("E21").Value = "  +0.33%  "

$ws.Range("D22").Value = "2.103.87"
$ws.Range("E22").Value = "  +0.39%  "

$ws.Range("E23").Value = "  +0.78%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.582"
$ws.Range("E24").Value = "  -0.44%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.27"
$ws.Range("E25").Value = "  -0.32%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.872"
$ws.Range("E26").Value = "  -1.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.170"
$ws.Range("E27").Value = "  +0.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.32"
$ws.Range("E28").Value = "  -0.28%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.61"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.866"
$ws.Range("E30").Value = "  -0.88%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08984"
$ws.Range("E31").Value = "  +0.75%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.868"
$ws.Range("E32").Value = "  +4.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.174"
$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("E34").Value = "  +0.77%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7742"
$ws.Range("E35").Value = "  +0.98%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.632"
$ws.Range("E36").Value = "  +4.09%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02087"
$ws.Range("E37").Value = "  +2.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.073"
$ws.Range("E38").Value = "  +3.00%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.095"
$ws.Range("E39").Value = "  +0.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5521"
$ws.Range("E40").Value = "  +0.83%  "

$ws.Range("E41").Value = "  +0.25%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.674"
$ws.Range("E42").Value = "  -3.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "114.54"
$ws.Range("E43").Value = "  +2.97%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.524"
$ws.Range("E44").Value = "  +0.56%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1512"
$ws.Range("E45").Value = "  -0.38%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4816"
$ws.Range("E46").Value = "  +0.42%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.47"
$ws.Range("E47").Value = "  -1.22%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9980"
$ws.Range("E48").Value = "  -0.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.618"
$ws.Range("E49").Value = "  -1.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "66.75"
$ws.Range("E50").Value = "  -0.98%  "

$ws.Range("E51").Value = "  -1.09%  "
